$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.383.02"
$ws.Range("E2").Value = "  -4.40%  "
$ws.Range("D3").Value = "1.567.58"
$ws.Range("E3").Value = "  -4.62%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("D7").Value = "'0.3668"
$ws.Range("E7").Value = "  -3.11%  "
$ws.Range("D8").Value = "'49.49"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("E9").Value = "  -3.37%  "
$ws.Range("D10").Value = "'1.171"
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("E11").Value = "  -5.73%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "'21.20"
$ws.Range("E13").Value = "  -3.82%  "
$ws.Range("E14").Value = "  -4.88%  "
$ws.Range("D15").Value = "'6.896"
$ws.Range("E15").Value = "  -5.67%  "
$ws.Range("D16").Value = "1.575.50"
$ws.Range("E16").Value = "  -4.07%  "
$ws.Range("D17").Value = "'0.00001135"
$ws.Range("E17").Value = "  -5.37%  "
$ws.Range("D18").Value = "'89.10"
$ws.Range("E18").Value = "  -7.61%  "
$ws.Range("D19").Value = "'0.06771"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D21").Value = "'6.247"
$ws.Range("E21").Value = "  -6.90%  "
$ws.Range("E22").Value = "  -6.27%  "
$ws.Range("E23").Value = "  -4.78%  "
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").Value = "22.388.23"
$ws.Range("E25").Value = "  -4.49%  "
$ws.Range("D26").Value = "'2.386"
$ws.Range("E26").Value = "  -4.17%  "
$ws.Range("D27").Value = "'2.975"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("D28").Value = "'19.89"
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("D30").Value = "'4.955"
$ws.Range("E30").Value = "  -4.79%  "
$ws.Range("D31").Value = "'125.30"
$ws.Range("E31").Value = "  -5.49%  "
$ws.Range("D32").Value = "1.751.42"
$ws.Range("E32").Value = "  -4.06%  "
$ws.Range("D33").Value = "'1.035"
$ws.Range("E33").Value = "  +5.76%  "
$ws.Range("D34").Value = "'6.250"
$ws.Range("E34").Value = "  -8.93%  "
$ws.Range("D35").Value = "'1.992"
$ws.Range("E35").Value = "  -5.95%  "
$ws.Range("D36").Value = "'10.30"
$ws.Range("D37").Value = "'0.08463"
$ws.Range("E37").Value = "  -3.09%  "
$ws.Range("D38").Value = "'0.02542"
$ws.Range("E38").Value = "  -5.78%  "
$ws.Range("D39").Value = "'0.2328"
$ws.Range("E39").Value = "  -4.09%  "
$ws.Range("D40").Value = "'0.06557"
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("D41").Value = "'5.532"
$ws.Range("E41").Value = "  -6.18%  "
$ws.Range("E42").Value = "  -7.89%  "
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("E44").Value = "  -7.07%  "
$ws.Range("D45").Value = "'14.49"
$ws.Range("E45").Value = "  -7.12%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'0.6012"
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("D49").Value = "'2.132"
$ws.Range("E49").Value = "  -5.12%  "
$ws.Range("D50").Value = "'1.273"
$ws.Range("E50").Value = "  +8.38%  "
$ws.Range("D51").Value = "'123.47"
$ws.Range("E51").Value = "  -2.78%  "
